$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion note text with refreshed rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 11.25 = 45782.81 pesos`n✅ 45782.81 pesos = 11.22 = 970.6 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- Update "tasas" sheet rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 88.898
$ws2.Range("O10").Value = 4070
$ws2.Range("N12").Value = 4081.1
$ws2.Range("O12").Value = 86.52
